# Applies the commit "Add files via upload" text edits to the STAMP SRPM deck.
#
# Strategy: for each target shape, locate the old substring inside the
# shape's full TextRange.Text and replace just that slice via
# TextRange.Characters(start, length) so the untouched characters keep
# their original run formatting (rPr) intact.

function Replace-Text {
    param($tr, [string]$old, [string]$new)

    $full = $tr.Text
    $idx = $full.IndexOf($old)
    if ($idx -lt 0) {
        Write-Host "NOT FOUND: [$old]"
        return
    }
    $sub = $tr.Characters($idx + 1, $old.Length)
    $sub.Text = $new
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 - "STAMP Requirements/Scope" bullets
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange

Replace-Text $tr3 "Eliminate per session provisioning" "Eliminate per session provisioning on Session-Reflector"
Replace-Text $tr3 "Stateless on session-reflector" "No control-channel signaling for sessions"
Replace-Text $tr3 "Support very high scale for number of sessions and faster detection interval" "Support hardware implementation - very high scale for number of sessions and faster detection interval"

# ---------------------------------------------------------------------
# Slide 4 - "Draft status" bullets
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange

Replace-Text $tr4 "Move Receive Counter and other Response message fields to Section 4.1 from 3.2" "Move Receive Counter and other Reply message fields to Section 4.1 from 3.2"

# ---------------------------------------------------------------------
# Slide 5 - "Session-Sender Control Code Field"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)

# Rectangle 2: drop the leading "In a Query: " run, and reword the
# Query/Response language to Test-packet/Reply language.
$sh5a = $s5.Shapes.Item(4)
$tr5a = $sh5a.TextFrame.TextRange

Replace-Text $tr5a "In a Query: " ""
Replace-Text $tr5a "0x0: Out-of-band Response Requested.  " "0x0: Out-of-band Reply Requested.  "
Replace-Text $tr5a "0x1: In-band Response Requested.  " "0x1: In-band Reply Requested.  "
Replace-Text $tr5a "Indicates that this query has been sent over a bidirectional path and the probe response is required over the same path in reverse direction." "Indicates that this test packet has been sent over a bidirectional path and the reply is required over the same path in reverse direction."
Replace-Text $tr5a "0x2: No Response Requested." "0x2: No Reply Requested."

# Rectangle 3: figure caption under the STAMP DM message diagram.
$sh5b = $s5.Shapes.Item(5)
$tr5b = $sh5b.TextFrame.TextRange

Replace-Text $tr5b "      Figure: Session-Sender Control Code in STAMP DM Message" "             Session-Sender Control Code in Test Packet"
# The caption used to be followed by a trailing space run; remove it now
# that the new caption text no longer needs it.
$full5b = $tr5b.Text
if ($full5b.Substring($full5b.Length - 1, 1) -eq " ") {
    $trailing = $tr5b.Characters($full5b.Length, 1)
    $trailing.Text = ""
}

# ---------------------------------------------------------------------
# Slide 7 - "Return Path TLV"
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(5)
$tr7 = $sh7.TextFrame.TextRange

Replace-Text $tr7 "Type (value 1): Return Address. Target node address of the response; different than the Source Address in the query" "Type (value 1): Return Address. Target node address for the reply; different than the Source Address in the test packet"

# ---------------------------------------------------------------------
# Slide 9 - "Destination Node Address TLV"
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(5)
$tr9 = $sh9.TextFrame.TextRange

Replace-Text $tr9 "Indicates the address of the intended recipient node of the query message.  " "Indicates the address of the intended recipient node of the test packet message.  "
Replace-Text $tr9 "send response if it is not the intended destination node of the query." "send reply if it is not the intended destination node of the test packet."
Replace-Text $tr9 "Useful when query is sent with 127/8 destination address (e.g. sweeping ECMP paths)." "Useful when test packet is sent with 127/8 destination address (e.g. sweeping ECMP paths)."

# The text grew by a few characters, which re-wraps this auto-fit text
# box onto an extra line; match PowerPoint's resulting box height.
$sh9.Height = 239.4912011023622

# ---------------------------------------------------------------------
# Slide 10 - "Stand-alone Direct-mode LM Message Format"
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(3)
$tr10 = $sh10.TextFrame.TextRange

Replace-Text $tr10 "Stand-alone Direct-mode Loss Measurement (LM) query and response messages defined" "Stand-alone Direct-mode Loss Measurement (LM) test packet defined"
Replace-Text $tr10 "Direct-mode LM message format is also defined for authenticated mode" "Direct-mode LM packet format is also defined for authenticated mode"
Replace-Text $tr10 " is used for identifying direct-mode LM probe packets" " is used for identifying direct-mode LM test packets"
